$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Helper: apply the bold / centered / thin-bordered "index" style used
# throughout this workbook for header cells and the leading index column.
function Set-IndexStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# ---------------------------------------------------------------------------
# Header row (row 1) of the new "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # starts at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    Set-IndexStyle $cell
}

# ---------------------------------------------------------------------------
# Data rows (rows 2..10) of the new "2022-Q1" sheet.
# Columns: A index(n) | B code(text) | C name(text) | D scale(text) |
#          E total-position(text) | F position-pct(text) | G value(text) | H rank(n)
# ---------------------------------------------------------------------------
$rows = @(
    @("011136", "广发盛兴混合A", "22.19", "93.21", "5.13", "1.1383", 8),
    @("010161", "广发瑞安精选股票A", "8.03", "92.30", "5.02", "0.4031", 7),
    @("011137", "广发盛兴混合C", "2.10", "93.21", "5.13", "0.1077", 8),
    @("010162", "广发瑞安精选股票C", "0.65", "92.30", "5.02", "0.0326", 7),
    @("006786", "泰康中证港股通大消费主题指数A", "0.85", "80.77", "3.57", "0.0303", 6),
    @("159735", "银华中证港股通消费主题交易型开放式指数证券投资基金", "0.72", "92.83", "3.71", "0.0267", 7),
    @("513230", "华夏中证港股通消费主题ETF", "0.60", "96.92", "3.86", "0.0232", 7),
    @("006787", "泰康中证港股通大消费主题指数C", "0.41", "80.77", "3.57", "0.0146", 6),
    @("513590", "鹏华中证港股通消费主题交易型开放式指数证券投资基金", "0.37", "91.21", "3.77", "0.0139", 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $i
    Set-IndexStyle $idxCell

    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $nameCell = $newSheet.Cells.Item($r, 3)
    $nameCell.NumberFormat = "@"
    $nameCell.Value = $row[1]

    $scaleCell = $newSheet.Cells.Item($r, 4)
    $scaleCell.NumberFormat = "@"
    $scaleCell.Value = $row[2]

    $posCell = $newSheet.Cells.Item($r, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[3]

    $pctCell = $newSheet.Cells.Item($r, 6)
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $row[4]

    $valCell = $newSheet.Cells.Item($r, 7)
    $valCell.NumberFormat = "@"
    $valCell.Value = $row[5]

    $rankCell = $newSheet.Cells.Item($r, 8)
    $rankCell.Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new first data row for 2022-Q1 and
#    renumber the leading index column for every existing data row.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()
$ws.Range("A2:D2").ClearFormats()

$idxCell = $ws.Cells.Item(2, 1)
$idxCell.Value = 0
Set-IndexStyle $idxCell
$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 9
$ws.Cells.Item(2, 4).Value = 1.79

# Renumber the A-column index for the rows that got shifted down (previously
# 0..4 in rows 2..6, now sitting in rows 3..7) so the sequence stays 1..5.
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
